$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range('B2').Style

$ws.Range('D2').Value = '29.290.46'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.873.14'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = $plainStyle
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7104'
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '242.58'
$ws.Range('D6').Style = $plainStyle
$ws.Range('E6').Value = '  +0.96%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = $plainStyle
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3109'
$ws.Range('D8').Style = $plainStyle
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07728'
$ws.Range('D9').Style = $plainStyle
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '25.00'
$ws.Range('D10').Style = $plainStyle
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08472'
$ws.Range('D11').Style = $plainStyle
$ws.Range('E11').Value = '  +2.70%  '
$ws.Range('D12').Value = '1.866.66'
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.198'
$ws.Range('D13').Style = $plainStyle
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7103'
$ws.Range('D14').Style = $plainStyle
$ws.Range('E14').Value = '  -0.81%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.21'
$ws.Range('D15').Style = $plainStyle
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.000008322'
$ws.Range('D16').Style = $plainStyle
$ws.Range('E16').Value = '  +6.80%  '
$ws.Range('D17').Value = '29.297.15'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.981'
$ws.Range('D18').Style = $plainStyle
$ws.Range('E18').Value = '  +2.11%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '242.42'
$ws.Range('D19').Style = $plainStyle
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').Value = '2.123.77'
$ws.Range('E21').Value = '  +0.80%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9997'
$ws.Range('D22').Style = $plainStyle
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.800'
$ws.Range('D23').Style = $plainStyle
$ws.Range('E23').Value = '  -1.89%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.000'
$ws.Range('D24').Style = $plainStyle
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1610'
$ws.Range('D25').Style = $plainStyle
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '162.97'
$ws.Range('D26').Style = $plainStyle
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.011'
$ws.Range('D27').Style = $plainStyle
$ws.Range('E27').Value = '  +1.16%  '
$ws.Range('E28').Value = '  +1.73%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.512'
$ws.Range('D29').Style = $plainStyle
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.407'
$ws.Range('D30').Style = $plainStyle
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.325'
$ws.Range('D31').Style = $plainStyle
$ws.Range('E31').Value = '  +5.62%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.263'
$ws.Range('D32').Style = $plainStyle
$ws.Range('E32').Value = '  -4.58%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05259'
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.916'
$ws.Range('D34').Style = $plainStyle
$ws.Range('E34').Value = '  +0.63%  '
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7438'
$ws.Range('D36').Style = $plainStyle
$ws.Range('E36').Value = '  +2.11%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.684'
$ws.Range('D37').Style = $plainStyle
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01858'
$ws.Range('D38').Style = $plainStyle
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.717'
$ws.Range('D39').Style = $plainStyle
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('D40').Value = '1.165.64'
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.351'
$ws.Range('D41').Style = $plainStyle
$ws.Range('E41').Value = '  +4.71%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '72.83'
$ws.Range('D42').Style = $plainStyle
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8859'
$ws.Range('D43').Style = $plainStyle
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '106.63'
$ws.Range('D44').Style = $plainStyle
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9999'
$ws.Range('D45').Style = $plainStyle
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '2.019.62'
$ws.Range('E46').Value = '  +0.66%  '
$ws.Range('E47').Value = '  +2.22%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5197'
$ws.Range('D48').Style = $plainStyle
$ws.Range('E48').Value = '  -1.35%  '
$ws.Range('E49').Value = '  +4.09%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.369'
$ws.Range('D50').Style = $plainStyle
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4297'
$ws.Range('D51').Style = $plainStyle
$ws.Range('E51').Value = '  +1.35%  '
